# Apply the edit described in the commit:
#   - fix header text from "isbn's" to "isbn"
#   - insert a duplicate ISBN row (9780201563177) above the existing one at row 3
#   - move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header cell text (A1): "isbn's" -> "isbn"
$ws.Range("A1").Value = "isbn"

# Insert a new row at row 3, shifting everything below down by one.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with a duplicate of the ISBN that is
# now on row 4 (the original row 3 content, pushed down).
$ws.Cells.Item(3, 1).Value = 9780201563177

# Update the current selection to match the saved workbook view.
$ws.Range("G15").Select()
